# Test_Cases.xlsx edit
# Strip the redundant "SC_" prefix from every Scenario ID in column A
# (e.g. "SC_LOGIN_01" -> "LOGIN_01") and update the view: zoom to 130%
# and move the active selection to A22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).SpecialCells(11).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = [string]$cell.Value2
    if ($current -like "SC_*") {
        $cell.Value2 = $current.Substring(3)
    }
}

$win = $excel.ActiveWindow
$win.Zoom = 130

$ws.Range("A22").Select()
